# Updated cryptos list on Mon Aug 26 05:31:11 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding price figures are plain text in the source data (they keep
# exact trailing/leading zeros, and some use "." as a thousands separator).
# Force Text format before writing so Excel does not silently reinterpret
# the literal string as a number and reformat/round it.
$ws.Range("D2").Value = "63.922.08"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.739.04"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.00"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.85"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.109"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "3.228.22"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.99"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "63.776.00"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000149"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").Value = "2.744.38"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.19"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.87"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.61"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  -5.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.32"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.42"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").Value = "0.0₃0908"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.97"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.23"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("E31").Value = "  +7.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.96"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.09"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.986"
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "347.95"
$ws.Range("E39").Value = "  +4.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.34"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.09"
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.62"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.86"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.07"
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0584"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.42"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.623"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.100"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0250"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.03%  "
